$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Sheet1: clear obsolete C/D columns, comments, and validations ----
$ws1.Range("C1:D1").Comment.Delete()
$ws1.Range("B2:B1048576").Validation.Delete()
$ws1.Range("C2:C1048576").Validation.Delete()
$ws1.Range("D2:D1048576").Validation.Delete()
$ws1.Range("C1:D1").Clear()

# ---- Sheet1: new header values + comments ----
$ws1.Range("A1").Value = "abc"
$ws1.Range("B1").Value = "xyz"
[void]$ws1.Range("A1").Comment.Text("Start of the alphabet")
[void]$ws1.Range("B1").Comment.Text("End of the alphabet")

# ---- Sheet1: new list validations referencing the enum values ----
$ws1.Range("A2:A1048576").Validation.Add(3, 1, 1, """A,B,C""")
$ws1.Range("B2:B1048576").Validation.Add(3, 1, 1, """X,Y,Z""")

# ---- Rename main sheet ----
$ws1.Name = "Export this as TSV"

# ---- Add "abc list" sheet right after sheet1 ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "abc list"
$ws2.Range("A1").Value = "A"
$ws2.Range("A2").Value = "B"
$ws2.Range("A3").Value = "C"

# ---- Add "xyz list" sheet right after "abc list" ----
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "xyz list"
$ws3.Range("A1").Value = "X"
$ws3.Range("A2").Value = "Y"
$ws3.Range("A3").Value = "Z"

# ---- Restore sheet1 as the active/selected sheet ----
$ws1.Activate()
